$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New results run (2021/04/20 19:45): updated score for person 5 / rater 2.
$ws.Range("B11").Value = 0.7608

$wb.Save()
